# Update "paises.xlsx" (Pais sheet) with refreshed COVID-19 country data.
# 1) Update the "last refreshed" timestamp banner in A1.
# 2) Push new Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes figures for the handful of
#    countries whose numbers changed in this refresh.
# 3) Re-sort the whole country table (A4:H219) descending by "Casos
#    totales" (column B), which is how this sheet is always ordered, so the
#    updated rows land in their correct new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 00:57"

# 2) Updated per-country figures: row -> (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Row numbers are the countries' positions *before* the re-sort below.
$ws.Range("B4").Value = 2885652
$ws.Range("C4").Value = 49968
$ws.Range("D4").Value = 1206912
$ws.Range("E4").Value = 1546697
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 558
$ws.Range("H4").Value = 132043

$ws.Range("B22").Value = 109505
$ws.Range("C22").Value = 3395
$ws.Range("D22").Value = 45334
$ws.Range("E22").Value = 60394
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 136
$ws.Range("H22").Value = 3777

$ws.Range("B28").Value = 72786
$ws.Range("C28").Value = 2845
$ws.Range("D28").Value = 25224
$ws.Range("E28").Value = 46125
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 52
$ws.Range("H28").Value = 1437

$ws.Range("B32").Value = 60657
$ws.Range("C32").Value = 1189
$ws.Range("D32").Value = 28391
$ws.Range("E32").Value = 27566
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 61
$ws.Range("H32").Value = 4700

$ws.Range("B53").Value = 27564
$ws.Range("C53").Value = 454
$ws.Range("D53").Value = 11069
$ws.Range("E53").Value = 15867
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 12
$ws.Range("H53").Value = 628

$ws.Range("B56").Value = 21293
$ws.Range("C56").Value = 2282
$ws.Range("D56").Value = 3315
$ws.Range("E56").Value = 17098
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 63
$ws.Range("H56").Value = 880

$ws.Range("B58").Value = 19068
$ws.Range("C58").Value = 194
$ws.Range("D58").Value = 16866
$ws.Range("E58").Value = 1226
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 976

$ws.Range("B69").Value = 12319
$ws.Range("C69").Value = 141
$ws.Range("D69").Value = 7848
$ws.Range("E69").Value = 4118
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 353

$ws.Range("B88").Value = 5620
$ws.Range("C88").Value = 107
$ws.Range("D88").Value = 2555
$ws.Range("E88").Value = 3021
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 44

$ws.Range("B89").Value = 5497
$ws.Range("C89").Value = 182
$ws.Range("D89").Value = 2892
$ws.Range("E89").Value = 2366
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 239

$ws.Range("B134").Value = 1082
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 963
$ws.Range("E134").Value = 51
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 68

$ws.Range("B168").Value = 277
$ws.Range("C168").Value = 50
$ws.Range("D168").Value = 29
$ws.Range("E168").Value = 247
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 1

$ws.Range("B176").Value = 146
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 137
$ws.Range("E176").Value = 0
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 9

# 3) Re-sort the country table descending by "Casos totales" (column B) so
#    the rows reflect the refreshed ranking (shared strings / row order in
#    the saved file follow this sort).
$sortRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$sortRange.Sort($sortKey, 2)
